$d = $word.ActiveDocument

# Locate the unique "doctor patient" text (from the Forthcoming paragraph)
# so we can turn it into "doctor-patient" while splitting the run exactly
# the way the target diff expects: "...doctor" / "-" / "patient..." as
# three separate runs instead of one run whose text merely changes.
$find = $d.Content.Find
$found = $find.Execute("doctor patient", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $matchRange = $find.Parent

    # Position of the single space between "doctor" and "patient".
    $spaceStart = $matchRange.Start + 6
    $spaceEnd = $spaceStart + 1

    # Replace the space with a hyphen.
    $hyphenRange = $d.Range($spaceStart, $spaceEnd)
    $hyphenRange.Text = "-"

    # Toggling a character property on the newly-typed hyphen forces Word
    # to keep it as its own run (rather than silently re-merging it back
    # into the neighbouring text once the formatting round-trips to its
    # original value), matching the three-run split produced by the
    # original edit: "...doctor" + "-" + "patient relationship...".
    $hyphenRange2 = $d.Range($spaceStart, $spaceStart + 1)
    $hyphenRange2.Bold = 1
    $hyphenRange2.Bold = 0
}
